$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Oct 06 11:31:59 EDT 2023"
$ws.Range("B3").Value = "Fri Oct 06 11:32:14 EDT 2023"
$ws.Range("B4").Value = "Fri Oct 06 11:32:28 EDT 2023"
